$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.066.52'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '1.567.22'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  +0.58%  '
$ws.Range('D5').Value = "'208.52"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.08%  '
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('D8').Value = "'22.08"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('D10').Value = "'0.0597"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('D11').Value = "'0.0861"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.565.41'
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'3.78"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = "'0.520"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.25%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '27.068.77'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = "'61.89"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.0₃0706'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = "'7.42"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.11%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = "'215.80"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = "'1.01"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = "'4.15"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = "'9.20"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').Value = "'1.94"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = "'154.06"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = "'6.62"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = "'15.04"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = "'0.105"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.29%  '
$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D28').Value = "'1.01"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = "'0.0475"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = "'1.13"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.95%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'3.24"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = "'3.20"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.12%  '
$ws.Range('B33').Value = 'Maker'
$ws.Range('C33').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D33').Value = '1.423.62'
$ws.Range('E33').Value = '  +0.77%  '
$ws.Range('B34').Value = 'TrustWalletToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D34').Value = "'1.09"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +12.58%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = "'1.61"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.35"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.0167"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.94%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = "'0.535"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = "'5.83"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.51%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = "'0.811"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = "'1.01"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.52%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = "'2.36"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.49%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = "'1.00"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = "'64.74"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = "'1.74"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.704.23'
$ws.Range('E46').Value = '  +1.11%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = "'86.80"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0103'
$ws.Range('E48').Value = '  +2.90%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = "'0.0518"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = "'0.0964"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = "'1.01"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.57%  '
